$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-45 down to 27-46.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly data point.
$ws.Cells.Item(26, 1).Value = 8
$ws.Cells.Item(26, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 45049
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 100112026
$ws.Cells.Item(26, 7).Value = "Haba"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 300
$ws.Cells.Item(26, 11).Value = 14000
$ws.Cells.Item(26, 12).Value = 15000
$ws.Cells.Item(26, 13).Value = 14500
$ws.Cells.Item(26, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(26, 16).Value = 580
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
